$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename column G header from "Status" to "Test Data"
$ws.Range("G1").Value = "Test Data"

# 2. Insert a new row at position 20 (pushes old rows 19.. down; row 19 "7A" stays in place)
$ws.Rows.Item(20).Insert()

# 3. Update row 19 ("7A" - Contact Us Form) with the full script content
$ws.Range("C19").Value = "Check user can submit an enquiry using valid entries"
$ws.Range("E19").Value = "1. Open https://abantecart.codifyme.co.nz/`n2. Click Contact Us on the footer block`n3. Enter First name`n4. Enter Email`n5. Enter Enquiry message`n6. Click Submit"
$ws.Range("F19").Value = "6. A message appears `"Your enquiry has been successfully sent to the store owner!`"`nhttps://abantecart.codifyme.co.nz/index.php?rt=content/contact/success"

# Row 19 now has complete content like the other fully fleshed-out test cases
# (rows 2-18), so copy their cell formatting (style) onto row 19.
$ws.Range("A2:F2").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Populate the newly inserted row 20 ("7B" - Contact Us Form)
$ws.Range("A20").Value = "7B"
$ws.Range("B20").Value = $ws.Range("B19").Value2
$ws.Range("C20").Value = "Check user cannot submit a blank enquiry"
$ws.Range("D20").Value = $ws.Range("D19").Value2
$ws.Range("E20").Value = "1. Open https://abantecart.codifyme.co.nz/`n2. Click Contact Us on the footer block`n3. Leave all fields blank and click Submit"
$ws.Range("F20").Value = "3. Each field returns with an error message:`n• First name: First name: is required field! Name must be between 3 and 32 characters!`n• Email: Email: is required field! E-Mail Address does not appear to be valid!`n• Enquiry: Enquiry: is required field! Enquiry must be between 10 and 3000 characters!"

# Match the row heights Excel's own word-wrap auto-fit produced for these two rows
$ws.Rows.Item(19).RowHeight = 90
$ws.Rows.Item(20).RowHeight = 105

# 5. Restore view state: frozen pane anchor + active selection
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A15").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("B19").Select()
